# Apply cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "28.737.25"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "1.574.52"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0890"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.799.83"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "1.577.54"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "28.741.48"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -4.74%  "
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "1.396.13"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.525"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.794"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0471"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.960"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "1.711.90"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +0.37%  "
